{"js": "// \"inducted clarity in experience\"\n// 1. \"...warranted employment...\" -> \"...warranted eventual employment...\"\n// 2. Remove the manual line break between \"focusing on \" and \"end-to-end\"\n//    so the sentence reads as one continuous line.\n\nconst body = context.document.body;\n\n// --- Edit 1: insert \"eventual \" before \"employment\" -----------------------\nconst warrantedResults = body.search(\"warranted employment\", { matchCase: true });\nwarrantedResults.load(\"text\");\nawait context.sync();\n\nif (warrantedResults.items.length > 0) {\n  warrantedResults.items[0].insertText(\n    \"warranted eventual employment\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// --- Edit 2: remove the manual line break before \"end-to-end\" -------------\n// The break character shows up in range text as \\u000b (vertical tab).\nconst breakResults = body.search(\"focusing on \\u000bend-to-end\", { matchCase: true });\nbreakResults.load(\"text\");\nawait context.sync();\n\nif (breakResults.items.length > 0) {\n  breakResults.items[0].insertText(\n    \"focusing on end-to-end\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n} else {\n  // Fallback: locate the break via the paragraph that contains \"end-to-end\"\n  // and rewrite its text without the embedded line-break character.\n  const paras = body.paragraphs;\n  paras.load(\"text\");\n  await context.sync();\n\n  for (let i = 0; i < paras.items.length; i++) {\n    const t = paras.items[i].text;\n    if (t.indexOf(\"\\u000b\") !== -1 && t.indexOf(\"end-to-end\") !== -1) {\n      paras.items[i].insertText(t.split(\"\\u000b\").join(\"\"), Word.InsertLocation.replace);\n      await context.sync();\n      break;\n    }\n  }\n}\n", "ps1": "# \"inducted clarity in experience\"\n# 1. \"...warranted employment...\" -> \"...warranted eventual employment...\"\n# 2. Remove the manual line break between \"focusing on \" and \"end-to-end\"\n#    so the sentence reads as one continuous line.\n\n$d = $word.ActiveDocument\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\n# --- Edit 1: insert \"eventual \" before \"employment\" ------------------------\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Text = \"warranted employment\"\n$find1.Replacement.Text = \"warranted eventual employment\"\n$found1 = $find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find1.Replacement.Text, $wdReplaceAll)\n\nif (-not $found1) {\n    # Fallback: locate the phrase's paragraph and run Find/Replace scoped to\n    # that paragraph's Range (avoids clobbering the paragraph mark that a\n    # direct Range.Text re-assignment would otherwise swallow).\n    foreach ($p in $d.Paragraphs) {\n        if ($p.Range.Text -like \"*warranted employment*\") {\n            $pf = $p.Range.Find\n            $pf.ClearFormatting()\n            $pf.Replacement.ClearFormatting()\n            $pf.Text = \"warranted employment\"\n            $pf.Replacement.Text = \"warranted eventual employment\"\n            $pf.Execute($pf.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $pf.Replacement.Text, $wdReplaceAll)\n            break\n        }\n    }\n}\n\n# --- Edit 2: remove the manual line break before \"end-to-end\" --------------\n# \"^l\" is Word's Find/Replace token for a manual line break (Chr(11)).\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"focusing on ^lend-to-end\"\n$find2.Replacement.Text = \"focusing on end-to-end\"\n$found2 = $find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find2.Replacement.Text, $wdReplaceAll)\n\nif (-not $found2) {\n    # Fallback: locate the paragraph containing the manual line break\n    # (Chr(11)) right before \"end-to-end\" and strip just that break via a\n    # Find/Replace scoped to the paragraph's Range.\n    foreach ($p in $d.Paragraphs) {\n        $t = $p.Range.Text\n        if ($t.IndexOf([char]11) -ge 0 -and $t -like \"*end-to-end*\") {\n            $pf = $p.Range.Find\n            $pf.ClearFormatting()\n            $pf.Replacement.ClearFormatting()\n            $pf.Text = \"^l\"\n            $pf.Replacement.Text = \"\"\n            $pf.Execute($pf.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $pf.Replacement.Text, $wdReplaceAll)\n            break\n        }\n    }\n}\n"}
